# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45170) {
        $cell.Value2 = 45174
    }
}
